$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (2-12, columns A-D) by column A ascending,
# matching the calibration data reordering.
$rng = $ws.Range("A2:D12")
$key = $ws.Range("A2:A12")
$rng.Sort($key, 1)
